$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '96.381.62'
$ws.Range("E2").Value = '  -1.08%  '
$ws.Range("D3").Value = '3.331.50'
$ws.Range("E3").Value = '  -1.98%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").Value = "'249.58"
$ws.Range("E5").Value = '  -1.94%  '
$ws.Range("D6").Value = "'653.76"
$ws.Range("E6").Value = '  +0.91%  '
$ws.Range("D7").Value = "'1.40"
$ws.Range("E7").Value = '  -3.37%  '
$ws.Range("D8").Value = "'0.420"
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").Value = "'0.993"
$ws.Range("E10").Value = '  -5.13%  '
$ws.Range("D11").Value = '3.328.80'
$ws.Range("E11").Value = '  -2.13%  '
$ws.Range("E12").Value = '  -2.56%  '
$ws.Range("D13").Value = "'40.29"
$ws.Range("E13").Value = '  -2.84%  '
$ws.Range("D14").Value = '96.121.94'
$ws.Range("E14").Value = '  -0.95%  '
$ws.Range("D15").Value = "'6.07"
$ws.Range("E15").Value = '  -2.70%  '
$ws.Range("E16").Value = '  -2.45%  '
$ws.Range("D17").Value = '3.941.60'
$ws.Range("E17").Value = '  -1.92%  '
$ws.Range("D18").Value = "'8.58"
$ws.Range("E18").Value = '  +1.59%  '
$ws.Range("D19").Value = '3.310.86'
$ws.Range("E19").Value = '  -2.45%  '
$ws.Range("B20").Value = 'Stellar'
$ws.Range("C20").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D20").Value = "'0.549"
$ws.Range("E20").Value = '  +12.89%  '
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").Value = "'17.14"
$ws.Range("E21").Value = '  -1.29%  '
$ws.Range("D22").Value = "'504.47"
$ws.Range("E22").Value = '  +0.48%  '
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").Value = "'10.51"
$ws.Range("E23").Value = '  -2.43%  '
$ws.Range("B24").Value = 'SuiNetwork'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D24").Value = "'3.36"
$ws.Range("E24").Value = '  -1.16%  '
$ws.Range("D25").Value = "'0.0000198"
$ws.Range("E25").Value = '  -1.96%  '
$ws.Range("D26").Value = "'6.58"
$ws.Range("E26").Value = '  +8.99%  '
$ws.Range("D27").Value = "'95.87"
$ws.Range("E27").Value = '  +0.23%  '
$ws.Range("D28").Value = "'12.03"
$ws.Range("E28").Value = '  -4.68%  '
$ws.Range("D29").Value = "'0.145"
$ws.Range("E29").Value = '  -5.29%  '
$ws.Range("E30").Value = '  +0.36%  '
$ws.Range("D31").Value = "'11.03"
$ws.Range("E31").Value = '  -2.13%  '
$ws.Range("D32").Value = "'0.188"
$ws.Range("E32").Value = '  -4.97%  '
$ws.Range("D33").Value = "'2.47"
$ws.Range("E33").Value = '  +10.27%  '
$ws.Range("E34").Value = '  +0.38%  '
$ws.Range("D35").Value = "'0.547"
$ws.Range("E35").Value = '  -3.23%  '
$ws.Range("D36").Value = "'28.00"
$ws.Range("E36").Value = '  -5.05%  '
$ws.Range("D37").Value = "'1.47"
$ws.Range("E37").Value = '  +5.81%  '
$ws.Range("D38").Value = "'7.63"
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("E39").Value = '  +0.03%  '
$ws.Range("E40").Value = '  -1.09%  '
$ws.Range("D41").Value = "'506.28"
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").Value = "'24.34"
$ws.Range("E42").Value = '  -1.45%  '
$ws.Range("D43").Value = "'0.0427"
$ws.Range("E43").Value = '  +4.15%  '
$ws.Range("D44").Value = "'0.830"
$ws.Range("E44").Value = '  -2.63%  '
$ws.Range("D45").Value = "'3.64"
$ws.Range("E45").Value = '  +0.09%  '
$ws.Range("D46").Value = "'1.67"
$ws.Range("E46").Value = '  +7.07%  '
$ws.Range("D47").Value = "'5.49"
$ws.Range("E47").Value = '  +1.11%  '
$ws.Range("D48").Value = "'8.36"
$ws.Range("E48").Value = '  +2.95%  '
$ws.Range("D49").Value = "'53.54"
$ws.Range("E49").Value = '  +4.17%  '
$ws.Range("D50").Value = "'3.11"
$ws.Range("E50").Value = '  -3.21%  '
$ws.Range("D51").Value = "'162.50"
$ws.Range("E51").Value = '  +1.44%  '
